{"js": "// Bold + enlarge (16pt) the five section-title paragraphs, and update the\n// \"Normal\" style's font color, matching the authored diff:\n//   - \"The Best of Reddit in 2017\"\n//   - \"Most Upvoted Posts of 2017\"\n//   - \"Top AMAs of 2017\"\n//   - \"Largest New Communities Created in 2017\"\n//   - \"Reddit Superlatives\"\n\nconst titles = [\n  \"The Best of Reddit in 2017\",\n  \"Most Upvoted Posts of 2017\",\n  \"Top AMAs of 2017\",\n  \"Largest New Communities Created in 2017\",\n  \"Reddit Superlatives\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text.trim();\n  if (titles.indexOf(text) !== -1) {\n    // Setting font on the paragraph (rather than on a sub-range) applies the\n    // formatting to both the run(s) of text and the paragraph mark itself.\n    para.font.bold = true;\n    para.font.size = 16;\n    // Also set the complex-script (bidirectional) counterparts so the bold\n    // size change is mirrored on w:bCs / w:szCs, matching Word's own\n    // behavior when bold/size is toggled from the UI.\n    para.font.boldBidirectional = true;\n    para.font.sizeBidirectional = 16;\n  }\n}\n\n// Update the \"Normal\" style's default font color (auto -> #00000A).\nconst styles = context.document.getStyles();\nstyles.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < styles.items.length; i++) {\n  styles.items[i].load(\"nameLocal\");\n}\nawait context.sync();\n\nconst normalStyle = styles.items.find((s) => s.nameLocal === \"Normal\");\nif (normalStyle) {\n  normalStyle.font.color = \"#00000A\";\n}\n\nawait context.sync();\n", "ps1": "# Bold + enlarge (16pt) the five section-title paragraphs, and update the\n# \"Normal\" style's font color, matching the authored diff:\n#   - \"The Best of Reddit in 2017\"\n#   - \"Most Upvoted Posts of 2017\"\n#   - \"Top AMAs of 2017\"\n#   - \"Largest New Communities Created in 2017\"\n#   - \"Reddit Superlatives\"\n\n$d = $word.ActiveDocument\n\n$titles = @(\n    \"The Best of Reddit in 2017\",\n    \"Most Upvoted Posts of 2017\",\n    \"Top AMAs of 2017\",\n    \"Largest New Communities Created in 2017\",\n    \"Reddit Superlatives\"\n)\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.Trim()\n    if ($titles -contains $text) {\n        # Use the paragraph's own Range (includes the paragraph mark) so the\n        # bold/size formatting lands on both the run and the paragraph mark,\n        # same as selecting the whole line (incl. pilcrow) in the Word UI.\n        $rng = $p.Range\n        $rng.Font.Bold = 1\n        $rng.Font.BoldBi = 1\n        $rng.Font.Size = 16\n        $rng.Font.SizeBi = 16\n    }\n}\n\n# Update the \"Normal\" style's default font color (auto -> #00000A == 655360).\n$normal = $d.Styles(\"Normal\")\n$normal.Font.Color = 655360\n"}
